$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells so numeric-looking strings
# (e.g. "1.004") are stored as literal text, not coerced to numbers.
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "24.353.54"
$ws.Range("E2").Value = "  +9.13%  "

$ws.Range("D3").Value = "1.676.60"
$ws.Range("E3").Value = "  +4.90%  "

$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").Value = "307.01"
$ws.Range("E5").Value = "  +6.23%  "

$ws.Range("D6").Value = "0.9982"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").Value = "0.3709"
$ws.Range("E7").Value = "  +0.40%  "

$ws.Range("D8").Value = "0.3436"
$ws.Range("E8").Value = "  +1.09%  "

$ws.Range("D9").Value = "48.11"
$ws.Range("E9").Value = "  +12.79%  "

$ws.Range("D10").Value = "1.179"
$ws.Range("E10").Value = "  +3.39%  "

$ws.Range("D11").Value = "0.07255"
$ws.Range("E11").Value = "  +2.95%  "

$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("D13").Value = "20.35"
$ws.Range("E13").Value = "  +3.20%  "

$ws.Range("D14").Value = "6.103"
$ws.Range("E14").Value = "  +2.95%  "

$ws.Range("D15").Value = "6.739"
$ws.Range("E15").Value = "  +1.24%  "

$ws.Range("D16").Value = "1.680.29"
$ws.Range("E16").Value = "  +5.14%  "

$ws.Range("D17").Value = "0.00001109"
$ws.Range("E17").Value = "  +2.31%  "

$ws.Range("D18").Value = "0.9986"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("D19").Value = "0.06727"
$ws.Range("E19").Value = "  +1.66%  "

$ws.Range("D20").Value = "81.04"
$ws.Range("E20").Value = "  +3.50%  "

$ws.Range("D21").Value = "16.43"
$ws.Range("E21").Value = "  +1.68%  "

$ws.Range("D22").Value = "6.097"
$ws.Range("E22").Value = "  +1.04%  "

$ws.Range("D23").Value = "11.94"
$ws.Range("E23").Value = "  +1.23%  "

$ws.Range("D24").Value = "24.295.03"
$ws.Range("E24").Value = "  +8.71%  "

$ws.Range("D25").Value = "2.435"
$ws.Range("E25").Value = "  +1.64%  "

$ws.Range("D26").Value = "3.365"
$ws.Range("E26").Value = "  -12.09%  "

$ws.Range("D27").Value = "2.655"
$ws.Range("E27").Value = "  +5.77%  "

$ws.Range("D28").Value = "152.51"
$ws.Range("E28").Value = "  +1.12%  "

$ws.Range("D29").Value = "19.55"
$ws.Range("E29").Value = "  -0.48%  "

$ws.Range("D30").Value = "1.865.07"
$ws.Range("E30").Value = "  +4.98%  "

$ws.Range("D31").Value = "127.18"
$ws.Range("E31").Value = "  +5.33%  "

$ws.Range("D32").Value = "6.314"
$ws.Range("E32").Value = "  +4.58%  "

$ws.Range("D33").Value = "4.026"
$ws.Range("E33").Value = "  -3.07%  "

$ws.Range("D34").Value = "0.9682"

$ws.Range("D35").Value = "1.738"
$ws.Range("E35").Value = "  +7.64%  "

$ws.Range("D36").Value = "0.08484"
$ws.Range("E36").Value = "  +2.39%  "

$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "12.33"
$ws.Range("E37").Value = "  +4.54%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.06455"
$ws.Range("E38").Value = "  +4.91%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "8.977"
$ws.Range("E39").Value = "  +3.57%  "

$ws.Range("D40").Value = "5.334"
$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("D41").Value = "0.02334"
$ws.Range("E41").Value = "  +5.32%  "

$ws.Range("D42").Value = "1.259"
$ws.Range("E42").Value = "  +1.67%  "

$ws.Range("D43").Value = "0.2107"
$ws.Range("E43").Value = "  +3.84%  "

$ws.Range("D44").Value = "0.6170"
$ws.Range("E44").Value = "  +4.01%  "

$ws.Range("D45").Value = "0.9983"
$ws.Range("E45").Value = "  +0.22%  "

$ws.Range("D46").Value = "3.776"
$ws.Range("E46").Value = "  +2.83%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.5941"
$ws.Range("E47").Value = "  +3.96%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "12.93"
$ws.Range("E48").Value = "  -1.74%  "

$ws.Range("D49").Value = "127.11"
$ws.Range("E49").Value = "  +0.99%  "

$ws.Range("D50").Value = "2.024"
$ws.Range("E50").Value = "  +2.59%  "

$ws.Range("D51").Value = "0.07202"
$ws.Range("E51").Value = "  +5.59%  "

# Restore default (General) styling on the Price cells so no
# stray style/number-format is left behind on the cells.
foreach ($addr in $dCells) {
    $ws.Range($addr).ClearFormats()
}
